# Update the "想去人数" (want-to-go count) values in column F across the
# four worksheets of the 广州-漫展信息 workbook, reflecting newly generated
# data (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) -> sheet1
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value  = 58
$wsExpo.Range("F3").Value  = 1019
$wsExpo.Range("F5").Value  = 14
$wsExpo.Range("F6").Value  = 457
$wsExpo.Range("F7").Value  = 740
$wsExpo.Range("F8").Value  = 255
$wsExpo.Range("F10").Value = 36
$wsExpo.Range("F11").Value = 409
$wsExpo.Range("F14").Value = 854
$wsExpo.Range("F16").Value = 1999
$wsExpo.Range("F17").Value = 493
$wsExpo.Range("F18").Value = 7559
$wsExpo.Range("F19").Value = 555
$wsExpo.Range("F24").Value = 226
$wsExpo.Range("F25").Value = 137

# Sheet "演出" (Performances) -> sheet2
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F8").Value  = 125
$wsShow.Range("F9").Value  = 57
$wsShow.Range("F10").Value = 7

# Sheet "本地生活" (Local Life) -> sheet3
$wsLocal = $wb.Worksheets.Item("本地生活")
$wsLocal.Range("F2").Value = 5535
$wsLocal.Range("F3").Value = 400
$wsLocal.Range("F4").Value = 391

# Sheet "全部类型" (All Types) -> sheet4, aggregates all of the above
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value  = 58
$wsAll.Range("F3").Value  = 5535
$wsAll.Range("F4").Value  = 400
$wsAll.Range("F5").Value  = 391
$wsAll.Range("F7").Value  = 1019
$wsAll.Range("F11").Value = 14
$wsAll.Range("F12").Value = 457
$wsAll.Range("F13").Value = 740
$wsAll.Range("F14").Value = 255
$wsAll.Range("F17").Value = 36
$wsAll.Range("F18").Value = 409
$wsAll.Range("F23").Value = 854
$wsAll.Range("F25").Value = 125
$wsAll.Range("F26").Value = 1999
$wsAll.Range("F27").Value = 493
$wsAll.Range("F28").Value = 7559
$wsAll.Range("F29").Value = 57
$wsAll.Range("F30").Value = 7
$wsAll.Range("F31").Value = 555
$wsAll.Range("F37").Value = 226
$wsAll.Range("F39").Value = 137
